# CDS Study filter fixes
# - Replace the "ParticipantsTab" Neo4j query (cell B2 on the "startup" sheet)
#   with the corrected query (adds diagnosis/file/genomic_info traversal and
#   sorts the collected sample ids).
# - Row 2 grows taller to fit the longer wrapped query text.
# - Selection moves to B2 (the cell that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newParticipantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Detection of Colorectal Cancer Susceptibility Loci Using Genome-Wide Sequencing"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

# Cell A2 ("ParticipantsTab") and the rest of the row are left untouched;
# only the query text in B2 is replaced.
$ws.Range("B2").Value = $newParticipantsQuery

# The wrapped text is now much longer, so the row needs to grow to show it.
$ws.Rows.Item(2).RowHeight = 299.25

# Move the active selection onto the cell that was edited.
$ws.Range("B2").Select() | Out-Null
